# Natmi following Dr Hou advice
# Rewrite the LR-pairs data block (rows 2-16) of Sheet1 with the recomputed
# NATMI statistics, which now include a third sending/receiving cluster
# ("M2") in addition to the original "FAPs" and "sCs" clusters.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs","Fgf7","Nrp1","ECs",3.0,1.0,18.32032066666667,54.960962,0.9393883922888365,0.9570930138987986,3.0,1.0,109.1447706666667,327.434312,0.3535542089399963,0.3655959674582361,1999.567197703127,17996.10477932814,0.3321247199230946,0.3499093463638503),
    @("FAPs","Fgf7","Nrp1","FAPs",3.0,1.0,18.32032066666667,54.960962,0.9393883922888365,0.9570930138987986,3.0,1.0,47.980825,143.942475,0.155424969272891,0.1607186127944892,879.0240998512168,7911.21689866095,0.1460044120068029,0.1538226615091117),
    @("FAPs","Fgf7","Nrp1","M1",3.0,1.0,18.32032066666667,54.960962,0.9393883922888365,0.9570930138987986,3.0,1.0,62.26741999999999,186.80226,0.2017037397794264,0.2085735992386923,1140.759101486013,10266.83191337412,0.1894781518300412,0.1996243347150802),
    @("FAPs","Fgf7","Nrp1","M2",3.0,1.0,18.32032066666667,54.960962,0.9393883922888365,0.9570930138987986,3.0,1.0,58.81030666666667,176.43092,0.1905050633580386,0.1969935053322898,1077.423676638338,9696.813089745041,0.1789582451907908,0.1885411077369704),
    @("FAPs","Fgf7","Nrp1","sCs",3.0,1.0,18.32032066666667,54.960962,0.9393883922888365,0.9570930138987986,2.0,1.0,30.503993,61.007986,0.09881201864964768,0.06811831517629259,558.8429333737554,3353.057600242532,0.09282286333810708,0.06519556357378616),
    @("M2","Fgf7","Nrp1","ECs",1.0,0.3333333333333333,0.09978633333333332,0.299359,0.005116620224500324,0.005213052994737072,3.0,1.0,109.1447706666667,327.434312,0.3535542089399963,0.3655959674582361,10.89115646733422,98.02040820600799,0.001809002615919599,0.001905871153021955),
    @("M2","Fgf7","Nrp1","FAPs",1.0,0.3333333333333333,0.09978633333333332,0.299359,0.005116620224500324,0.005213052994737072,3.0,1.0,47.980825,143.942475,0.155424969272891,0.1607186127944892,4.787830597058333,43.090475373525,0.0007952505411740155,0.0008378346457382997),
    @("M2","Fgf7","Nrp1","M1",1.0,0.3333333333333333,0.09978633333333332,0.299359,0.005116620224500324,0.005213052994737072,3.0,1.0,62.26741999999999,186.80226,0.2017037397794264,0.2085735992386923,6.213437527926666,55.92093775133999,0.001032041434312764,0.001087305226134355),
    @("M2","Fgf7","Nrp1","M2",1.0,0.3333333333333333,0.09978633333333332,0.299359,0.005116620224500324,0.005213052994737072,3.0,1.0,58.81030666666667,176.43092,0.1905050633580386,0.1969935053322898,5.868464864475555,52.81618378028,0.000974742060047456,0.001026937582916247),
    @("M2","Fgf7","Nrp1","sCs",1.0,0.3333333333333333,0.09978633333333332,0.299359,0.005116620224500324,0.005213052994737072,2.0,1.0,30.503993,61.007986,0.09881201864964768,0.06811831517629259,3.043881613495667,18.263289680974,0.0005055835730464906,0.0003551043869262159),
    @("sCs","Fgf7","Nrp1","ECs",2.0,1.0,1.082285,2.16457,0.05549498748666317,0.03769393310646423,3.0,1.0,109.1447706666667,327.434312,0.3535542089399963,0.3655959674582361,118.1257481209733,708.7544887258399,0.01962048640098219,0.01378074994136383),
    @("sCs","Fgf7","Nrp1","FAPs",2.0,1.0,1.082285,2.16457,0.05549498748666317,0.03769393310646423,3.0,1.0,47.980825,143.942475,0.155424969272891,0.1607186127944892,51.928927185125,311.57356311075,0.008625306724914093,0.006058116639639202),
    @("sCs","Fgf7","Nrp1","M1",2.0,1.0,1.082285,2.16457,0.05549498748666317,0.03769393310646423,3.0,1.0,62.26741999999999,186.80226,0.2017037397794264,0.2085735992386923,67.39109465469998,404.3465679281999,0.01119354651507243,0.007861959297477749),
    @("sCs","Fgf7","Nrp1","M2",2.0,1.0,1.082285,2.16457,0.05549498748666317,0.03769393310646423,3.0,1.0,58.81030666666667,176.43092,0.1905050633580386,0.1969935053322898,63.64951275073333,381.8970765044,0.01057207610720033,0.007425460012403238),
    @("sCs","Fgf7","Nrp1","sCs",2.0,1.0,1.082285,2.16457,0.05549498748666317,0.03769393310646423,2.0,1.0,30.503993,61.007986,0.09881201864964768,0.06811831517629259,33.014014064005,132.05605625602,0.005483571738494126,0.00256764721558022)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

Write-Output "Rows written"
